$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value  = -22.36520000000001
$ws.Range("A10").Value = -21.9879
$ws.Range("A12").Value = -21.5469
$ws.Range("A18").Value = -22.29790000000002
$ws.Range("A37").Value = -19.60429999999999
$ws.Range("A55").Value = -22.36280000000001
$ws.Range("A68").Value = -21.54799999999999
$ws.Range("A77").Value = -20.80469999999999
$ws.Range("A78").Value = -20.47659999999998
$ws.Range("A81").Value = -21.8202
$ws.Range("A82").Value = -21.9856
